# Tff2-Cxcr4.xlsx: refresh NATMI LR-pair edge table with new TPM-based values.
# A new "ECs" sending cluster is introduced (rows 2-6), and the existing
# Inflammatory-Mac / Resolving-Mac blocks are recomputed and a new
# Resolving-Mac target row set is appended, growing the table from 10 to 15 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = @(
    "Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity"
)

$dataRows = @(
    @("ECs", "Tff2", "Cxcr4", "ECs", 2, 0.6666666666666666, 0.5224156666666667, 1.567247, 0.1417347608275029, 0.1417347608275029, 3, 1, 18.10884, 54.32652, 0.02835750805894763, 0.02862140385105863, 9.460341721159999, 85.14307549044, 0.004019244622398928, 0.004056647829377166),
    @("ECs", "Tff2", "Cxcr4", "FAPs", 2, 0.6666666666666666, 0.5224156666666667, 1.567247, 0.1417347608275029, 0.1417347608275029, 1, 0.3333333333333333, 0.4738613333333333, 1.421584, 0.0007420423715060531, 0.0007489478393278889, 0.2475525843608889, 2.227973259248, 0.0001051731980492835, 0.0001061519428794134),
    @("ECs", "Tff2", "Cxcr4", "Inflammatory-Mac", 2, 0.6666666666666666, 0.5224156666666667, 1.567247, 0.1417347608275029, 0.1417347608275029, 3, 1, 305.722738, 917.168214, 0.4787460161246407, 0.4832012404300546, 159.7143479874287, 1437.429131886858, 0.06785495209254579, 0.06848641224390653),
    @("ECs", "Tff2", "Cxcr4", "MuSCs", 2, 0.6666666666666666, 0.5224156666666667, 1.567247, 0.1417347608275029, 0.1417347608275029, 2, 1, 17.6638495, 35.327699, 0.02766067592116823, 0.01861205798213543, 9.227871712442166, 55.367230274653, 0.003920479286013848, 0.002637975586605582),
    @("ECs", "Tff2", "Cxcr4", "Resolving-Mac", 2, 0.6666666666666666, 0.5224156666666667, 1.567247, 0.1417347608275029, 0.1417347608275029, 3, 1, 296.6213786666667, 889.864136, 0.4644937575237376, 0.4688163498974233, 154.9596552837324, 1394.636897553592, 0.06583491162849506, 0.06644757322473421),
    @("Inflammatory-Mac", "Tff2", "Cxcr4", "ECs", 3, 1, 1.746494666666667, 5.239484, 0.4738353377607538, 0.4738353377607539, 3, 1, 18.10884, 54.32652, 0.02835750805894763, 0.02862140385105863, 31.62699247951999, 284.64293231568, 0.01343678940916475, 0.01356183256095331),
    @("Inflammatory-Mac", "Tff2", "Cxcr4", "FAPs", 3, 1, 1.746494666666667, 5.239484, 0.4738353377607538, 0.4738353377607539, 1, 0.3333333333333333, 0.4738613333333333, 1.421584, 0.0007420423715060531, 0.0007489478393278889, 0.8275962914062222, 7.448366622656, 0.0003516058977353614, 0.000354877952413117),
    @("Inflammatory-Mac", "Tff2", "Cxcr4", "Inflammatory-Mac", 3, 1, 1.746494666666667, 5.239484, 0.4738353377607538, 0.4738353377607539, 3, 1, 305.722738, 917.168214, 0.4787460161246407, 0.4832012404300546, 533.9431313957307, 4805.488182561577, 0.2268467802520344, 0.2289578229655902),
    @("Inflammatory-Mac", "Tff2", "Cxcr4", "MuSCs", 3, 1, 1.746494666666667, 5.239484, 0.4738353377607538, 0.4738353377607539, 2, 1, 17.6638495, 35.327699, 0.02766067592116823, 0.01861205798213543, 30.84981894455266, 185.098913667316, 0.0131066057177975, 0.008819050780387878),
    @("Inflammatory-Mac", "Tff2", "Cxcr4", "Resolving-Mac", 3, 1, 1.746494666666667, 5.239484, 0.4738353377607538, 0.4738353377607539, 3, 1, 296.6213786666667, 889.864136, 0.4644937575237376, 0.4688163498974233, 518.0476558606471, 4662.428902745824, 0.2200935564840219, 0.2221417535014094),
    @("Resolving-Mac", "Tff2", "Cxcr4", "ECs", 3, 1, 1.416958, 4.250874, 0.3844299014117433, 0.3844299014117433, 3, 1, 18.10884, 54.32652, 0.02835750805894763, 0.02862140385105863, 25.65946570872, 230.93519137848, 0.01090147402738395, 0.01100292346072816),
    @("Resolving-Mac", "Tff2", "Cxcr4", "FAPs", 3, 1, 1.416958, 4.250874, 0.3844299014117433, 0.3844299014117433, 1, 0.3333333333333333, 0.4738613333333333, 1.421584, 0.0007420423715060531, 0.0007489478393278889, 0.6714416071573333, 6.042974464415999, 0.0002852632757214082, 0.0002879179440353585),
    @("Resolving-Mac", "Tff2", "Cxcr4", "Inflammatory-Mac", 3, 1, 1.416958, 4.250874, 0.3844299014117433, 0.3844299014117433, 3, 1, 305.722738, 917.168214, 0.4787460161246407, 0.4832012404300546, 433.196279391004, 3898.766514519036, 0.1840442837800605, 0.185757005220558),
    @("Resolving-Mac", "Tff2", "Cxcr4", "MuSCs", 3, 1, 1.416958, 4.250874, 0.3844299014117433, 0.3844299014117433, 2, 1, 17.6638495, 35.327699, 0.02766067592116823, 0.01861205798213543, 25.02893285982099, 150.173597158926, 0.01063359091735688, 0.007155031615141975),
    @("Resolving-Mac", "Tff2", "Cxcr4", "Resolving-Mac", 3, 1, 1.416958, 4.250874, 0.3844299014117433, 0.3844299014117433, 3, 1, 296.6213786666667, 889.864136, 0.4644937575237376, 0.4688163498974233, 420.3000354727627, 3782.700319254864, 0.1785652894112206, 0.1802270231712798)
)

for ($col = 0; $col -lt $headerRow.Length; $col++) {
    $ws.Cells.Item(1, $col + 1).Value = $headerRow[$col]
}

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $dataRows[$i]
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($rowNum, $col + 1).Value = $rowValues[$col]
    }
}

